# Auto-generated edit script: update cryptos price/volume columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.733.99"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "1.638.87"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'217.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").Value = "'0.250"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("D9").Value = "'0.0623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "'19.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.25%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.868.17"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "1.629.84"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").Value = "'4.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").Value = "'64.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.85%  "
$ws.Range("D17").Value = "26.729.13"
$ws.Range("E17").Value = "  +0.26%  "
$ws.Range("D18").Value = "0.0₃0728"
$ws.Range("E18").Value = "  -2.08%  "
$ws.Range("D19").Value = "'211.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.17%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "'4.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("E23").Value = "  +3.60%  "
$ws.Range("D24").Value = "'9.23"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "'145.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").Value = "'0.117"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.89%  "
$ws.Range("D28").Value = "'7.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("D29").Value = "'15.60"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").Value = "'0.0505"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'3.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").Value = "1.276.70"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  -0.72%  "
$ws.Range("D36").Value = "'2.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").Value = "'0.0175"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").Value = "'0.532"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("D39").Value = "'0.813"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("D41").Value = "'0.804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("E42").Value = "  -1.78%  "
$ws.Range("E43").Value = "  -0.20%  "
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("D45").Value = "'60.79"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("D46").Value = "'91.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("D48").Value = "'0.0522"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").Value = "'7.54"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("E51").Value = "  +0.09%  "
